# This script applies a weekly update to the "Femacal de La Calera - Kiwi" sheet.
# A new week of price data (fecha = 2021-10-05, serial 44474) is inserted as three
# new rows right before the existing block of data (which starts at row 377),
# pushing all the subsequent rows down by three. The worksheet dimension grows
# from A1:T433 to A1:T436 automatically as a result.
#
# NOTE: this runtime's PowerShell engine only supports *positional* parameter
# binding for user-defined functions (named "-Param value" binding does not
# work), so Set-KiwiRow below is called positionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the top of the data block (rows 377-379).
# Everything that used to live at row 377 onward shifts down to row 380 onward.
$ws.Rows("377:379").Insert()

function Set-KiwiRow {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $RegionDetalle, $PrecioUnit, $Kilos)

    $ws.Cells.Item($Row, 1).Value  = 3
    $ws.Cells.Item($Row, 2).Value  = "Femacal de La Calera"
    $ws.Cells.Item($Row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = 5
    $ws.Cells.Item($Row, 6).Value  = "Fruta"
    $ws.Cells.Item($Row, 7).Value  = 100101
    $ws.Cells.Item($Row, 8).Value  = "Berries"
    $ws.Cells.Item($Row, 9).Value  = 100101007
    $ws.Cells.Item($Row, 10).Value = "Kiwi"
    $ws.Cells.Item($Row, 11).Value = "Hayward"
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $RegionDetalle
    $ws.Cells.Item($Row, 19).Value = $PrecioUnit
    $ws.Cells.Item($Row, 20).Value = $Kilos
}

# New week of data: fecha 2021-10-05 (serial 44474).
Set-KiwiRow 377 44474 "Especial" 58 12000 12000 12000 "`$/bandeja 10 kilos" "Región de O'Higgins" 1200 10
Set-KiwiRow 378 44474 "Primera"  57 11000 11000 11000 "`$/bandeja 10 kilos" "Región de O'Higgins" 1100 10
Set-KiwiRow 379 44474 "Segunda"  50 10000 10000 10000 "`$/bandeja 10 kilos" "Región de O'Higgins" 1000 10
